$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update 想去人数 (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2047
$wsExhibit.Range("F5").Value = 1133

# Sheet "全部类型" (all types) - same underlying rows, update mirrored values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2047
$wsAll.Range("F7").Value = 1133
